# Upload new version with timestamp
# A new stock item ("قطن 100 جم") is inserted as row 34 (item #28), pushing the
# previously-last item ("كريم ONE للبشره الجافه") down to row 35 (renumbered
# #29), the grand-total row down to row 36 (total increased by the new
# item's sell price), and the footer row down to row 37 (with an updated
# timestamp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown = -4121
$xlPasteFormats = -4122

# 1) Insert a fresh blank row above the current row 34. This shifts the
#    existing rows 34 (cream item), 35 (grand total) and 36 (footer) down to
#    35, 36 and 37 respectively, including their merged-cell ranges.
$ws.Range("A34:Q34").Insert($xlShiftDown)

# 2) Clone the formatting (styles/number formats/fonts) of the row that is
#    now in position 35 (the original row 34) onto the newly inserted row 34
#    so the new item looks identical to the rest of the table.
$ws.Range("A35:Q35").Copy()
$ws.Range("A34:Q34").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Re-create the merged cells for the new row (PasteSpecial does not carry
# merge state over) and fix up the row heights to match the final layout.
$ws.Range("A34:B34").Merge()
$ws.Range("C34:G34").Merge()
$ws.Range("H34:K34").Merge()
$ws.Range("L34:M34").Merge()
$ws.Range("N34:O34").Merge()

$ws.Rows.Item(34).RowHeight = 25.5
$ws.Rows.Item(35).RowHeight = 24.75
$ws.Rows.Item(36).RowHeight = 25.5

# 3) Populate the new row 34 with the new item's data. L34/P34 keep
#    numeric-looking number formats inherited from the pasted formatting
#    (like every other row), but the underlying values are plain text, so
#    each is written while the cell is temporarily set to a text format -
#    this stores a real text value without Excel re-interpreting "0" /
#    "40.0000" as numbers, then the original numeric-looking format is
#    restored.
$ws.Range("A34").Value = 28
$ws.Range("C34").Value = "قطن 100 جم"
$ws.Range("H34").Value = "19:0"

$l34Format = $ws.Range("L34").NumberFormat
$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value = "0"
$ws.Range("L34").NumberFormat = $l34Format

$ws.Range("N34").Value = "20.00"

$p34Format = $ws.Range("P34").NumberFormat
$ws.Range("P34").NumberFormat = "@"
$ws.Range("P34").Value = "40.0000"
$ws.Range("P34").NumberFormat = $p34Format

$ws.Range("Q34").Value = "2:0"

# 4) The old item #28 is now row 35 - renumber it to #29.
$ws.Range("A35").Value = 29

# 5) Update the grand-total row (now row 36) with the new total.
$ws.Range("P36").Value = 1830.54

# 6) Update the footer row (now row 37) with the refreshed timestamp.
$ws.Range("A37").Value = "Tuesday, 19 August, 2025 12:24 PM"
